$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.933.44'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.910.75'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3810'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07304'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9110'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.27'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07663'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.914.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.491'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.67'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008735'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9978'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '27.985.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.192'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.128.94'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.87'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.614'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.45'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.837'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.202'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.43'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '115.37'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.916'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09024'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.203'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.843'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.234'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7786'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.649'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02088'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.093'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.096'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5547'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.05271'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.842'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '114.63'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.531'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4825'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9993'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.637'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '67.39'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06056'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.34%  '
